$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: add Reviewer-2 Chapter-2 structural response (B16) and thesis-change note (C16) ---
$ws.Range("B16").Value2 = "Both reviewers have suggested structural change to Chapter 2. I have followed the suggestions of Tim Butler as they slightly reduce the length of the chapter and may partly mitigate the issues mentioned here. Although I have moved Figure 2.37 as suggested"

$ws.Range("C16").Value2 = "I have implemented Tim Butler’s suggestions (first response for second reviewer in this document). and TODO: Fig 2.37 has been moved earlier, with references updated"
$ws.Range("C16").Font.Name = "Ubuntu"
$ws.Range("C16").Font.Size = 6.4
$ws.Range("C16").Font.Color = 3947580
$ws.Range("C16").WrapText = $true

$ws.Rows.Item(16).RowHeight = 208.45

# --- Row 28: add response to Reviewer-2's general Chapter-2 comment (B28), and extend C28 ---
$ws.Range("B28").Value2 = "Thank you for the constructive feedback, I have done my best to handle each item mentioned here."
$ws.Range("C28").Value2 = "Some intro from chapter 2 has been removed or replaced into chapter 1, with references added in chapter 2. Section 2.7: Filtering Data has been moved into Chapter 3 as subsection 3.2.3 under methods. Contributions and acknowledgements has been added as a small section (2.8) prior to Data Access, listing summarily who did what for the chapter and where some work has been used outside the thesis."

$ws.Rows.Item(28).RowHeight = 361.9

# --- Row 34: extend C34 with a TODO marker ---
$ws.Range("C34").Value2 = "Chemistry is not my strongest suit, so I do not pick up errors as readily as I should – thanks for pointing these out. It appears that ozone is photolysed by light up to long wavelengths (~1100nm) however the cross section is orders of magnitude lower for wavelengths past 320nm (TODO: reference)."

$ws.Rows.Item(34).RowHeight = 361.9

# --- cosmetic: scroll/zoom/selection state as left by the editing session ---
$ws.Application.ActiveWindow.Zoom = 100
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Range("C16").Select()

# --- column width tweaks (minor re-flow caused by the new content) ---
$ws.Columns.Item(1).ColumnWidth = 87.1122448979592
$ws.Columns.Item(2).ColumnWidth = 56.9642857142857
$ws.Columns.Item(3).ColumnWidth = 69.7448979591837
